$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 426
$ws.Range("I9").Value = 573.46155
$ws.Range("J9").Value = 106.5
$ws.Range("K9").Value = 573.46155
$ws.Range("L9").Value = 106.5
$ws.Range("M9").Value = -404.46155
$ws.Range("N9").Value = -444.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 19000
$ws.Range("I31").Value = 19000
$ws.Range("K31").Value = 57000
$ws.Range("M31").Value = -56770

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 1243.4
$ws.Range("I48").Value = 108.5
$ws.Range("K48").Value = 325.5
$ws.Range("M48").Value = -33.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H56").Value = 1243.4
$ws.Range("I56").Value = 108.5
$ws.Range("K56").Value = 325.5
$ws.Range("M56").Value = 208.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2755.2222

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2755.2222

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 83333976
$ws.Range("I98").Value = 100000570
$ws.Range("K98").Value = 100000570
$ws.Range("M98").Value = -99999072

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 83333976
$ws.Range("I122").Value = 100000570
$ws.Range("K122").Value = 300001710
$ws.Range("M122").Value = -299999260

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2645.5103
$ws.Range("I137").Value = 651.8929000000001
$ws.Range("J137").Value = 5303.6665
$ws.Range("K137").Value = 1955.6787
$ws.Range("L137").Value = 15910.9995
$ws.Range("M137").Value = 594.3212999999998
$ws.Range("N137").Value = -21010.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1565.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6956660
$ws.Range("I32").Value = 7943066.5
$ws.Range("K32").Value = 7943066.5
$ws.Range("M32").Value = -7942779.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 31328890
$ws.Range("I61").Value = 100012200
$ws.Range("J61").Value = 109203.27
$ws.Range("K61").Value = 100012200
$ws.Range("L61").Value = 109203.27
$ws.Range("M61").Value = -100011988
$ws.Range("N61").Value = -109627.27

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5686329
$ws.Range("I74").Value = 7143888
$ws.Range("K74").Value = 7143888
$ws.Range("M74").Value = -7143014

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5686329
$ws.Range("I77").Value = 7143888
$ws.Range("K77").Value = 35719440
$ws.Range("M77").Value = -35715072

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1565.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 31328890
$ws.Range("I136").Value = 100012200
$ws.Range("J136").Value = 109203.27
$ws.Range("K136").Value = 300036600
$ws.Range("L136").Value = 327609.81
$ws.Range("M136").Value = -300034050
$ws.Range("N136").Value = -332709.81

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1565.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3192.25
$ws.Range("J86").Value = 3089.6667
$ws.Range("L86").Value = 3089.6667
$ws.Range("N86").Value = -5335.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3192.25
$ws.Range("J89").Value = 3089.6667
$ws.Range("L89").Value = 15448.3335
$ws.Range("N89").Value = -26680.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 61939
$ws.Range("I96").Value = 14999
$ws.Range("J96").Value = 108879
$ws.Range("K96").Value = 14999
$ws.Range("L96").Value = 108879
$ws.Range("M96").Value = -12253
$ws.Range("N96").Value = -114371

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 32787.03
$ws.Range("I134").Value = 2515.3667
$ws.Range("K134").Value = 7546.1001
$ws.Range("M134").Value = -5011.1001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 694917
$ws.Range("I31").Value = 17696.334
$ws.Range("J31").Value = 938716.4
$ws.Range("K31").Value = 17696.334
$ws.Range("L31").Value = 938716.4
$ws.Range("M31").Value = -17401.334
$ws.Range("N31").Value = -939306.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 694917
$ws.Range("I34").Value = 17696.334
$ws.Range("J34").Value = 938716.4
$ws.Range("K34").Value = 17696.334
$ws.Range("L34").Value = 938716.4
$ws.Range("M34").Value = -17494.334
$ws.Range("N34").Value = -939120.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1014.6667
$ws.Range("I58").Value = 686.6667
$ws.Range("J58").Value = 1342.6666
$ws.Range("K58").Value = 686.6667
$ws.Range("L58").Value = 1342.6666
$ws.Range("M58").Value = -483.6667
$ws.Range("N58").Value = -1748.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 180999.8
$ws.Range("J74").Value = 180999.8
$ws.Range("L74").Value = 180999.8
$ws.Range("N74").Value = -182747.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 180999.8
$ws.Range("J77").Value = 180999.8
$ws.Range("L77").Value = 542999.3999999999
$ws.Range("N77").Value = -551735.3999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 837132.8
$ws.Range("I134").Value = 837132.8
$ws.Range("K134").Value = 2511398.4
$ws.Range("M134").Value = -2508863.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 95999.5
$ws.Range("J135").Value = 95999.5
$ws.Range("L135").Value = 95999.5
$ws.Range("N135").Value = -106139.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1014.6667
$ws.Range("I136").Value = 686.6667
$ws.Range("J136").Value = 1342.6666
$ws.Range("K136").Value = 2060.0001
$ws.Range("L136").Value = 4027.9998
$ws.Range("M136").Value = 489.9998999999998
$ws.Range("N136").Value = -9127.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 149
$ws.Range("I10").Value = 160.75
$ws.Range("K10").Value = 482.25
$ws.Range("M10").Value = -343.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 694.25
$ws.Range("J122").Value = 751
$ws.Range("L122").Value = 6759
$ws.Range("N122").Value = -11659

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 15172.167
$ws.Range("J125").Value = 15172.167
$ws.Range("L125").Value = 45516.501
$ws.Range("N125").Value = -55356.501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1550.3572
$ws.Range("J131").Value = 1648.25
$ws.Range("L131").Value = 4944.75
$ws.Range("N131").Value = -15024.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 249999.33
$ws.Range("J133").Value = 249999.33
$ws.Range("L133").Value = 249999.33
$ws.Range("N133").Value = -260119.33

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 201962.6
$ws.Range("I7").Value = 1267.6666
$ws.Range("J7").Value = 503005
$ws.Range("K7").Value = 1267.6666
$ws.Range("L7").Value = 503005
$ws.Range("M7").Value = -1155.6666
$ws.Range("N7").Value = -503229

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 25666.334
$ws.Range("I48").Value = 25666.334
$ws.Range("K48").Value = 25666.334
$ws.Range("M48").Value = -25005.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 125001010
$ws.Range("I55").Value = 166667710
$ws.Range("J55").Value = 876
$ws.Range("K55").Value = 166667710
$ws.Range("L55").Value = 876
$ws.Range("M55").Value = -166667537
$ws.Range("N55").Value = -1222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2804
$ws.Range("I82").Value = 1706
$ws.Range("K82").Value = 1706
$ws.Range("M82").Value = -1345

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2804
$ws.Range("I85").Value = 1706
$ws.Range("K85").Value = 1706
$ws.Range("M85").Value = -458

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 10092.546
$ws.Range("I100").Value = 11314.8
$ws.Range("J100").Value = 9074
$ws.Range("K100").Value = 11314.8
$ws.Range("L100").Value = 9074
$ws.Range("M100").Value = -10773.8
$ws.Range("N100").Value = -10156

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 201962.6
$ws.Range("I126").Value = 1267.6666
$ws.Range("J126").Value = 503005
$ws.Range("K126").Value = 3802.9998
$ws.Range("L126").Value = 1509015
$ws.Range("M126").Value = -1332.9998
$ws.Range("N126").Value = -1513955

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 91681.586
$ws.Range("I136").Value = 55492.316
$ws.Range("K136").Value = 166476.948
$ws.Range("M136").Value = -163926.948

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 9999.75
$ws.Range("J14").Value = 9999.75
$ws.Range("L14").Value = 9999.75
$ws.Range("N14").Value = -10335.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 50000530
$ws.Range("J107").Value = 560.6
$ws.Range("L107").Value = 1681.8
$ws.Range("N107").Value = -5521.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 359656.66
$ws.Range("I132").Value = 2322.3462
$ws.Range("J132").Value = 5005002.5
$ws.Range("K132").Value = 6967.0386
$ws.Range("L132").Value = 15015007.5
$ws.Range("M132").Value = -4437.0386
$ws.Range("N132").Value = -15020067.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3556.3333
$ws.Range("I136").Value = 939.0833
$ws.Range("K136").Value = 2817.2499
$ws.Range("M136").Value = -267.2498999999998
